# New weekly price data for "Cilantro" arrived (week of date serial 44523).
# It belongs at the top of the date-descending block that starts at row 104,
# so insert two fresh rows there (shifting the existing rows 104:145 down to
# 106:147, which is why dimension grows from A1:R145 to A1:R147) and fill
# them in with the same Primera/Segunda pair layout used by every other
# week in this block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("104:105").Insert()

$ws.Range("A104").Value = 11
$ws.Range("B104").Value = "Vega Monumental Concepción"
$ws.Range("C104").Value = "Bíobío"
$ws.Range("D104").Value = 44523
$ws.Range("E104").Value = 8
$ws.Range("F104").Value = 100112040
$ws.Range("G104").Value = "Cilantro"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 200
$ws.Range("K104").Value = 600
$ws.Range("L104").Value = 700
$ws.Range("M104").Value = 650
$ws.Range("N104").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O104").Value = "Región de Ñuble"
$ws.Range("P104").Value = 650
$ws.Range("Q104").Value = 1
$ws.Range("R104").Value = "Hortaliza"

$ws.Range("A105").Value = 11
$ws.Range("B105").Value = "Vega Monumental Concepción"
$ws.Range("C105").Value = "Bíobío"
$ws.Range("D105").Value = 44523
$ws.Range("E105").Value = 8
$ws.Range("F105").Value = 100112040
$ws.Range("G105").Value = "Cilantro"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Segunda"
$ws.Range("J105").Value = 100
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 500
$ws.Range("M105").Value = 500
$ws.Range("N105").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O105").Value = "Región de Ñuble"
$ws.Range("P105").Value = 500
$ws.Range("Q105").Value = 1
$ws.Range("R105").Value = "Hortaliza"
